# Applies the 05-11-2023 scraper update to the Primera RFEF Group 2
# 2023-2024 results sheet:
#   1) A batch of rows had their match data (columns F:V — home team,
#      score, away team, score, odds + timestamps, URL) re-ordered while
#      keeping the row's Indice/pais/torneio/temporada/data_partida
#      (columns A:E) fixed in place. Some are simple 2-row swaps, two are
#      3-row rotations.
#   2) Four new fixtures (rows 99-102) were appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Two-way swaps of F:V (match result/odds data) between row pairs ---
$pairs = @(
    @(12,13),
    @(14,15),
    @(23,24),
    @(29,30),
    @(35,37),
    @(38,39),
    @(53,55),
    @(57,58),
    @(85,86),
    @(93,94)
)

foreach ($pair in $pairs) {
    $a = $pair[0]
    $b = $pair[1]
    $rangeA = $ws.Range("F" + $a + ":V" + $a)
    $rangeB = $ws.Range("F" + $b + ":V" + $b)
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

# --- Three-way rotations of F:V between row groups ---
# new25 = old26, new26 = old28, new28 = old25
$v25 = $ws.Range("F25:V25").Value2
$v26 = $ws.Range("F26:V26").Value2
$v28 = $ws.Range("F28:V28").Value2
$ws.Range("F25:V25").Value = $v26
$ws.Range("F26:V26").Value = $v28
$ws.Range("F28:V28").Value = $v25

# new64 = old66, new65 = old64, new66 = old65
$v64 = $ws.Range("F64:V64").Value2
$v65 = $ws.Range("F65:V65").Value2
$v66 = $ws.Range("F66:V66").Value2
$ws.Range("F64:V64").Value = $v66
$ws.Range("F65:V65").Value = $v64
$ws.Range("F66:V66").Value = $v65

# --- Append 4 new fixtures (rows 99-102) ---
$newRows = @(
    @{ Row=99;  A=98;  E=45234.66666666666; F="Antequera";         G=5; H="Granada CF B"; I=2;
       J=1.68; K="02/11/2023 08:13"; L=1.78; M="04/11/2023 13:03";
       N=3.38; O="02/11/2023 08:13"; P=3.34; Q="04/11/2023 14:01";
       R=4.7;  S="02/11/2023 08:13"; T=5.02; U="04/11/2023 13:03";
       V="https://www.betexplorer.com/football/spain/primera-rfef-group-2/antequera-granada-cf/bqUmnWh5/" },
    @{ Row=100; A=99;  E=45234.75;          F="Recreativo Huelva"; G=1; H="San Fernando";  I=0;
       J=2.26; K="02/11/2023 08:13"; L=2.48; M="04/11/2023 17:47";
       N=2.95; O="02/11/2023 08:13"; P=2.88; Q="04/11/2023 16:05";
       R=3.12; S="02/11/2023 08:13"; T=3.27; U="04/11/2023 17:47";
       V="https://www.betexplorer.com/football/spain/primera-rfef-group-2/recreativo-huelva-cd-san-fernando/KhVqmCwa/" },
    @{ Row=101; A=100; E=45234.83333333334; F="Linares";           G=1; H="Algeciras";     I=1;
       J=2.75; K="02/11/2023 08:13"; L=3.01; M="04/11/2023 19:52";
       N=2.8;  O="02/11/2023 08:13"; P=2.98; Q="04/11/2023 19:52";
       R=2.63; S="02/11/2023 08:13"; T=2.56; U="04/11/2023 19:52";
       V="https://www.betexplorer.com/football/spain/primera-rfef-group-2/linares-algeciras/nJu4rlhU/" },
    @{ Row=102; A=101; E=45234.83333333334; F="CF Intercity";      G=1; H="Sanluqueno";    I=0;
       J=2.43; K="02/11/2023 08:13"; L=2.17; M="04/11/2023 19:54";
       N=3.13; O="02/11/2023 08:13"; P=3.15; Q="04/11/2023 19:54";
       R=2.77; S="02/11/2023 08:13"; T=3.6;  U="04/11/2023 19:54";
       V="https://www.betexplorer.com/football/spain/primera-rfef-group-2/cf-intercity-sanluqueno/QVVulhOh/" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $nr.A
    $cellA.Font.Bold = $true
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4160
    $cellA.Borders.LineStyle = 1

    $ws.Cells.Item($r, 2).Value = "spain"
    $ws.Cells.Item($r, 3).Value = "primera-rfef-group-2"
    $ws.Cells.Item($r, 4).Value = "2023-2024"

    $cellE = $ws.Cells.Item($r, 5)
    $cellE.Value = $nr.E
    $cellE.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 6).Value  = $nr.F
    $ws.Cells.Item($r, 7).Value  = $nr.G
    $ws.Cells.Item($r, 8).Value  = $nr.H
    $ws.Cells.Item($r, 9).Value  = $nr.I
    $ws.Cells.Item($r, 10).Value = $nr.J
    $ws.Cells.Item($r, 11).Value = $nr.K
    $ws.Cells.Item($r, 12).Value = $nr.L
    $ws.Cells.Item($r, 13).Value = $nr.M
    $ws.Cells.Item($r, 14).Value = $nr.N
    $ws.Cells.Item($r, 15).Value = $nr.O
    $ws.Cells.Item($r, 16).Value = $nr.P
    $ws.Cells.Item($r, 17).Value = $nr.Q
    $ws.Cells.Item($r, 18).Value = $nr.R
    $ws.Cells.Item($r, 19).Value = $nr.S
    $ws.Cells.Item($r, 20).Value = $nr.T
    $ws.Cells.Item($r, 21).Value = $nr.U
    $ws.Cells.Item($r, 22).Value = $nr.V
}
